# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-price / profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ - columns H:N)
# for the leves whose prices moved since the last run, across the ALC, ARM,
# BSM, CRP, CUL, GSM and LTW sheets. Rows/columns not listed are left
# untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 2332.6667
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 2749
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 8247
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -8583

# row 28
$ws.Range("H28").Value = 453
$ws.Range("I28").Value = 453
$ws.Range("K28").Value = 453
$ws.Range("M28").Value = 32

# row 40
$ws.Range("H40").Value = 7063.269
$ws.Range("J40").Value = 7575.9565
$ws.Range("L40").Value = 7575.9565
$ws.Range("N40").Value = -7925.9565

# row 41
$ws.Range("H41").Value = 1333
$ws.Range("I41").Value = 2650
$ws.Range("J41").Value = 674.5
$ws.Range("K41").Value = 2650
$ws.Range("L41").Value = 674.5
$ws.Range("M41").Value = -2210
$ws.Range("N41").Value = -1554.5

# row 98
$ws.Range("H98").Value = 1857.2667
$ws.Range("I98").Value = 1575.25
$ws.Range("K98").Value = 1575.25
$ws.Range("M98").Value = -77.25

# row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# row 122
$ws.Range("H122").Value = 1857.2667
$ws.Range("I122").Value = 1575.25
$ws.Range("K122").Value = 4725.75
$ws.Range("M122").Value = -2275.75

# row 132
$ws.Range("H132").Value = 4333.8823
$ws.Range("I132").Value = 3619.2856
$ws.Range("K132").Value = 10857.8568
$ws.Range("M132").Value = -8327.856800000001

$ws = $wb.Worksheets.Item("ARM")
# row 31
$ws.Range("H31").Value = 30236.5
$ws.Range("I31").Value = 30236.5
$ws.Range("K31").Value = 30236.5
$ws.Range("M31").Value = -29942.5

# row 107
$ws.Range("H107").Value = 99995
$ws.Range("J107").Value = 99995
$ws.Range("L107").Value = 99995
$ws.Range("N107").Value = -107675

# row 110
$ws.Range("H110").Value = 863
$ws.Range("I110").Value = 880.38464
$ws.Range("K110").Value = 880.38464
$ws.Range("M110").Value = 1164.61536

# row 121
$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489

# row 122
$ws.Range("H122").Value = 1704.6666
$ws.Range("I122").Value = 1625.6
$ws.Range("K122").Value = 4876.799999999999
$ws.Range("M122").Value = -2426.799999999999

# row 132
$ws.Range("H132").Value = 2855.5
$ws.Range("I132").Value = 2711
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8133
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5603
$ws.Range("N132").Value = -14060

# row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 670
$ws.Range("I22").Value = 615
$ws.Range("K22").Value = 615
$ws.Range("M22").Value = -442

# row 107
$ws.Range("H107").Value = 1037.7778
$ws.Range("I107").Value = 942.5
$ws.Range("K107").Value = 942.5
$ws.Range("M107").Value = 977.5

# row 108
$ws.Range("H108").Value = 95092
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()

# row 115
$ws.Range("H115").Value = 99995
$ws.Range("J115").Value = 99995
$ws.Range("L115").Value = 99995
$ws.Range("N115").Value = -103129

$ws = $wb.Worksheets.Item("CRP")
# row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# row 30
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# row 31
$ws.Range("H31").Value = 4075.6
$ws.Range("J31").Value = 7478.1665
$ws.Range("L31").Value = 7478.1665
$ws.Range("N31").Value = -8068.1665

# row 34
$ws.Range("H34").Value = 4075.6
$ws.Range("J34").Value = 7478.1665
$ws.Range("L34").Value = 7478.1665
$ws.Range("N34").Value = -7882.1665

# row 58
$ws.Range("H58").Value = 3386.4443
$ws.Range("J58").Value = 2166.3333
$ws.Range("L58").Value = 2166.3333
$ws.Range("N58").Value = -2572.3333

# row 107
$ws.Range("H107").Value = 1649.75
$ws.Range("I107").Value = 1599.6666
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 1599.6666
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 320.3334
$ws.Range("N107").Value = -5640

# row 122
$ws.Range("H122").Value = 2122.375
$ws.Range("I122").Value = 1530.6666
$ws.Range("J122").Value = 3897.5
$ws.Range("K122").Value = 4591.9998
$ws.Range("L122").Value = 11692.5
$ws.Range("M122").Value = -2141.9998
$ws.Range("N122").Value = -16592.5

# row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# row 132
$ws.Range("H132").Value = 3952.75
$ws.Range("I132").Value = 3952.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11858.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9328.25
$ws.Range("N132").ClearContents()

# row 134
$ws.Range("H134").Value = 957.875
$ws.Range("I134").Value = 957.875
$ws.Range("K134").Value = 2873.625
$ws.Range("M134").Value = -338.625

# row 136
$ws.Range("H136").Value = 3386.4443
$ws.Range("J136").Value = 2166.3333
$ws.Range("L136").Value = 6498.999899999999
$ws.Range("N136").Value = -11598.9999

# row 141
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360

$ws = $wb.Worksheets.Item("CUL")
# row 70
$ws.Range("H70").Value = 6000
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630

# row 73
$ws.Range("H73").Value = 6000
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184

# row 117
$ws.Range("H117").Value = 635.5
$ws.Range("J117").Value = 767
$ws.Range("L117").Value = 2301
$ws.Range("N117").Value = -9185

# row 127
$ws.Range("H127").Value = 666833
$ws.Range("J127").Value = 666833
$ws.Range("L127").Value = 2000499
$ws.Range("N127").Value = -2010419

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 1962.5
$ws.Range("I102").Value = 1962.5
$ws.Range("K102").Value = 1962.5
$ws.Range("M102").Value = -340.5

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 6824.125
$ws.Range("I40").Value = 6824.125
$ws.Range("K40").Value = 6824.125
$ws.Range("M40").Value = -6688.125

# row 122
$ws.Range("H122").Value = 5978.2
$ws.Range("I122").Value = 5923.25
$ws.Range("K122").Value = 17769.75
$ws.Range("M122").Value = -15319.75

